# Week1 L1 Exercise1: add a "Distinction?" column (J) to the Student
# Assignment sheet that flags students who passed all three tests, and
# switch the active sheet/selection back to the Student Assignment tab.

$wb = $excel.ActiveWorkbook
$wsStudent = $wb.Worksheets.Item("Student Assignment")

# --- New header cell J4: "Distinction?", styled like the neighbouring
#     Pass/Fail header cells (I4). ---
$wsStudent.Cells.Item(4, 10).Value = "Distinction?"
$wsStudent.Cells.Item(4, 9).Copy()
$wsStudent.Cells.Item(4, 10).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- New formula column J5:J12, mirrors the style used by G:I ---
for ($r = 5; $r -le 12; $r++) {
    $formula = '=IF(AND(G' + $r + '="Pass",H' + $r + '="Pass",I' + $r + '="Pass"),"Distinction","")'
    $wsStudent.Cells.Item($r, 10).Formula = $formula
}
$wsStudent.Cells.Item(5, 9).Copy()
$wsStudent.Range("J5:J12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column J sizing to fit the new header ---
$wsStudent.Columns.Item(10).ColumnWidth = 12

# --- Row heights grew slightly once the taller header wrapped ---
$wsStudent.Rows.Item(3).RowHeight = 15.75
$wsStudent.Rows.Item(4).RowHeight = 60.75

# --- Selection / active sheet: work moved back onto Student Assignment ---
$wsStudent.Activate()
$wsStudent.Range("H15").Select() | Out-Null
